$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2458139657974243
$ws.Range("B1").Value = 1.072569966316223
$ws.Range("C1").Value = 3.871024608612061
$ws.Range("D1").Value = 1.614328861236572
$ws.Range("E1").Value = 1.11142110824585
